# StudyTracker.xlsx — "adding some more files"
#
# Adds a new line item to the "OnLinuxPlatform" sheet (cell A5): a note
# about studying Android Linux Kernel internals once the base Linux OS
# material is done. Writing this value appends a new shared string and
# grows the sheet's used range/dimension to A1:A5.
#
# The newly-typed row keeps the same wrapped-text look as the row above
# it (row 4), so its height is set to match. Finally, the OnLinuxPlatform
# tab (the first sheet) is made the active tab/selection, replacing DBMS
# (previously tabSelected) as the sheet shown when the workbook opens.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("OnLinuxPlatform")

$ws.Range("A5").Value = "Once Linux OS is completed concentrate on Android Linux Kernel functionalities along with Android Architecture."

# Match row 4's wrapped-text row height for the newly added row.
$ws.Rows.Item(5).RowHeight = 30

# Make OnLinuxPlatform the active/selected sheet (was DBMS), with A5 --
# the cell just filled in -- as the active selection.
$ws.Select() | Out-Null
$ws.Range("A5").Select() | Out-Null
